$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextCell "D2" "50.850.88"
$ws.Range("E2").Value = "  -2.76%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.731.63"
$ws.Range("E3").Value = "  -3.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
Set-TextCell "D5" "349.71"
$ws.Range("E5").Value = "  -3.18%  "

# Row 6 - Solana
Set-TextCell "D6" "106.36"
$ws.Range("E6").Value = "  -4.38%  "

# Row 7 - XRP
Set-TextCell "D7" "0.543"
$ws.Range("E7").Value = "  -3.97%  "

# Row 8 - USDC
Set-TextCell "D8" "0.999"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.572"
$ws.Range("E9").Value = "  -4.60%  "

# Row 10 - Avalanche
Set-TextCell "D10" "38.67"
$ws.Range("E10").Value = "  -4.61%  "

# Row 11 - TRON
Set-TextCell "D11" "0.135"
$ws.Range("E11").Value = "  +2.51%  "

# Row 12 - Dogecoin
Set-TextCell "D12" "0.0825"
$ws.Range("E12").Value = "  -4.13%  "

# Row 13 - Chainlink
Set-TextCell "D13" "19.41"
$ws.Range("E13").Value = "  -2.05%  "

# Row 14 - Polkadot
Set-TextCell "D14" "7.39"
$ws.Range("E14").Value = "  -4.79%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell "D15" "3.174.04"
$ws.Range("E15").Value = "  -2.62%  "

# Row 16 - WrappedEther
Set-TextCell "D16" "2.781.20"
$ws.Range("E16").Value = "  -1.94%  "

# Row 17 - Polygon
Set-TextCell "D17" "0.909"
$ws.Range("E17").Value = "  -1.28%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "50.788.60"
$ws.Range("E18").Value = "  -2.44%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +0.85%  "

# Row 20 - ImmutableX
Set-TextCell "D20" "3.00"
$ws.Range("E20").Value = "  -4.34%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextCell "D21" "12.81"
$ws.Range("E21").Value = "  -3.90%  "

# Row 22 - ShibaInu
Set-TextCell "D22" "0.0₃0947"
$ws.Range("E22").Value = "  -4.77%  "

# Row 23 - Litecoin
Set-TextCell "D23" "68.76"
$ws.Range("E23").Value = "  -1.92%  "

# Row 24 - BitcoinCash
Set-TextCell "D24" "261.12"
$ws.Range("E24").Value = "  -4.52%  "

# Row 25 - PancakeSwap
Set-TextCell "D25" "2.68"
$ws.Range("E25").Value = "  -4.61%  "

# Row 26 - Dai
Set-TextCell "D26" "0.999"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "25.63"
$ws.Range("E27").Value = "  -4.38%  "

# Row 28 - Kaspa
Set-TextCell "D28" "0.159"
$ws.Range("E28").Value = "  +12.03%  "

# Row 29 - Toncoin
Set-TextCell "D29" "2.24"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30 - Cosmos
Set-TextCell "D30" "9.93"
$ws.Range("E30").Value = "  -3.12%  "

# Row 31 - OKB
Set-TextCell "D31" "51.26"
$ws.Range("E31").Value = "  -1.84%  "

# Row 32 - was Filecoin, now InjectiveProtocol
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D32" "34.12"
$ws.Range("E32").Value = "  -0.83%  "

# Row 33 - was InjectiveProtocol, now Filecoin
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D33" "5.95"
$ws.Range("E33").Value = "  +1.96%  "

# Row 34 - VeChain
Set-TextCell "D34" "0.0437"
$ws.Range("E34").Value = "  -8.37%  "

# Row 35 - Hedera
Set-TextCell "D35" "0.0820"
$ws.Range("E35").Value = "  -3.36%  "

# Row 36 - RenderToken
Set-TextCell "D36" "5.09"
$ws.Range("E36").Value = "  -7.15%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.59%  "

# Row 38 - Celestia
Set-TextCell "D38" "18.15"
$ws.Range("E38").Value = "  +0.18%  "

# Row 39 - LidoDAOToken
Set-TextCell "D39" "3.10"
$ws.Range("E39").Value = "  -3.77%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -5.38%  "

# Row 41 - Stellar
Set-TextCell "D41" "0.112"
$ws.Range("E41").Value = "  -3.79%  "

# Row 42 - Stacks
Set-TextCell "D42" "2.44"
$ws.Range("E42").Value = "  -4.16%  "

# Row 43 - WEMIXToken
Set-TextCell "D43" "2.20"
$ws.Range("E43").Value = "  -3.09%  "

# Row 44 - Monero
Set-TextCell "D44" "119.04"
$ws.Range("E44").Value = "  -4.84%  "

# Row 45 - EnergySwap
Set-TextCell "D45" "21.32"
$ws.Range("E45").Value = "  -4.90%  "

# Row 46 - Maker
Set-TextCell "D46" "2.057.38"
$ws.Range("E46").Value = "  -0.64%  "

# Row 47 - ApeXProtocol
Set-TextCell "D47" "2.31"
$ws.Range("E47").Value = "  -0.98%  "

# Row 48 - NEARProtocol
Set-TextCell "D48" "3.16"
$ws.Range("E48").Value = "  -4.02%  "

# Row 49 - was SEI, now THORChain
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell "D49" "5.40"
$ws.Range("E49").Value = "  -7.67%  "

# Row 50 - was THORChain, now SEI
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
Set-TextCell "D50" "0.902"
$ws.Range("E50").Value = "  -5.05%  "

# Row 51 - MultiversX
Set-TextCell "D51" "57.49"
$ws.Range("E51").Value = "  -4.24%  "
